# "refactor: wrap table recognizer"
# The long single-column table (rows 1-21, columns A-C) is reflowed into a
# two-block-per-row layout: the first 11 records stay in columns A-C, and
# records 12-21 are moved up into columns D-F alongside rows 2-11, with a
# repeated header in D1:F1. Some reference-value cells that had no data are
# cleared out (still present but blank) rather than carrying a value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for the second block (D1:F1) ---
$ws.Range("D1").Value = "项目名称"
$ws.Range("E1").Value = "结果"
$ws.Range("F1").Value = "参考值"

# --- Clear reference-value cells in column C that no longer carry a value ---
$ws.Range("C2").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("C7").Value = ""
$ws.Range("C8").Value = ""

# --- Move the former rows 12-21 (A/B, and occasionally C) up into D:F,
#     aligned with rows 2-11 next to the first block. The "result" column
#     keeps its original text formatting (e.g. "53.0", "1.10") instead of
#     being auto-converted to a number. ---
$ws.Range("D2").Value = "CREA肌酐"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "53.0"
$ws.Range("F2").Value = ""

$ws.Range("D3").Value = "尿酸"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "294"
$ws.Range("F3").Value = ""

$ws.Range("D4").Value = "葡萄糖"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.68"
$ws.Range("F4").Value = "3.8-6.2"

$ws.Range("D5").Value = "总胆固醇"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.74"
$ws.Range("F5").Value = "2.9-6.0"

$ws.Range("D6").Value = "甘油三酯"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.72"
$ws.Range("F6").Value = "0.55-1.7"

$ws.Range("D7").Value = "载脂蛋白A1"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.1"
$ws.Range("F7").Value = ""

$ws.Range("D8").Value = "载脂蛋白A1"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.8"
$ws.Range("F8").Value = ""

$ws.Range("D9").Value = "载脂蛋白比值"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.38"
$ws.Range("F9").Value = ""

$ws.Range("D10").Value = "乙肝表面抗原"
$ws.Range("E10").Value = "阴性"
$ws.Range("F10").Value = ""

$ws.Range("D11").Value = "B微球蛋白"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.15"
$ws.Range("F11").Value = ""

# --- The old rows 12-21 are no longer part of the table; remove them ---
$ws.Range("A12:C21").Clear()
